# Daily attendance processing - 2025-10-22 01:19:26
# Normalize the "Recorded By" (column G) strings: move a lone trailing
# "System" token to the front of the comma-separated recorder list, and
# for the triple "backup@backdoor.com, system, System" swap the last two
# tokens so the canonical-cased "System" precedes the lowercase duplicate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End($xlUp).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$replacements = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($null -eq $current) { continue }

    if ($replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
